$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New open/close/high/low price values per row (D, E, F, G)
$values = @{
    2  = @(36.49796724875108, 39.3697624206543, 39.3697624206543, 36.23081262504466)
    3  = @(39.17663253636952, 38.30202484130859, 39.17663253636952, 37.487412325024)
    4  = @(39.96397228402768, 40.53743743896485, 41.12443754708946, 38.71809526174923)
    5  = @(34.80086284580503, 33.37158584594727, 35.61070115027061, 32.67307002839149)
    6  = @(40.21736579158045, 38.85866546630859, 40.21736579158045, 38.20419139100507)
    7  = @(42.93101673829687, 42.29515838623047, 43.00026770789886, 41.25335384901548)
    8  = @(44.08157169323555, 45.70761489868164, 46.08213531289125, 43.79935209006902)
    9  = @(48.59069480363591, 46.50393676757812, 48.94718943000097, 46.48512187778989)
    10 = @(43.63729858398438, 44.44160079956055, 45.96900177001953, 42.59999847412109)
    11 = @(41.68000030517578, 43.81499862670898, 43.81999969482422, 37.27999877929688)
    12 = @(48.95240020751953, 49.81549835205078, 51.80500030517578, 48.79999923706055)
}

foreach ($r in $values.Keys) {
    $v = $values[$r]
    $ws.Range("D$r").Value = $v[0]
    $ws.Range("E$r").Value = $v[1]
    $ws.Range("F$r").Value = $v[2]
    $ws.Range("G$r").Value = $v[3]

    # shares_outstanding column is no longer populated
    $ws.Range("H$r").ClearContents()

    # fixed_ticker becomes "TEMP" for every row
    $ws.Range("I$r").Value = "TEMP"
}
